$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 7.9000000000000181
$ws.Range("C9").Value = 8.1033333333333513
$ws.Range("C10").Value = 7.8266666666666485
$ws.Range("C11").Value = 6.9199999999999484
$ws.Range("C12").Value = 6.2066666666666936
$ws.Range("C13").Value = 6.266666666666687
$ws.Range("C14").Value = 6.2200000000000033
$ws.Range("C16").Value = 6.843333333333379
$ws.Range("C17").Value = 6.9166666666666377
$ws.Range("C18").Value = 6.6633333333333544
$ws.Range("C19").Value = 7.1566666666666556
$ws.Range("C20").Value = 7.9833333333333423
$ws.Range("C21").Value = 8.470000000000022
$ws.Range("C22").Value = 9.4433333333333813
$ws.Range("C23").Value = 9.7266666666666612
$ws.Range("C24").Value = 9.0833333333332877
$ws.Range("C25").Value = 8.6133333333333404
$ws.Range("C26").Value = 8.2500000000000462
$ws.Range("C27").Value = 8.2433333333333358
$ws.Range("C28").Value = 8.1600000000000339
$ws.Range("C29").Value = 7.7433333333333465
$ws.Range("C30").Value = 6.4266666666667138
$ws.Range("C31").Value = 5.8633333333333093
$ws.Range("C32").Value = 5.6433333333333113
$ws.Range("C33").Value = 4.8166666666666691
$ws.Range("C34").Value = 4.0233333333333565
$ws.Range("C35").Value = 3.7700000000000289
$ws.Range("C36").Value = 3.2566666666666855
$ws.Range("C37").Value = 3.0366666666666653
$ws.Range("C38").Value = 3.0399999999999983
$ws.Range("C81").Value = 0.99666666666664572
$ws.Range("C82").Value = 1.0033333333333116
$ws.Range("C83").Value = 1.0095615055630391
$ws.Range("C84").Value = 1.4326571222309292
$ws.Range("C85").Value = 1.949196345653248
$ws.Range("C86").Value = 2.4691399677738479
$ws.Range("C87").Value = 2.9424672522763151
$ws.Range("C88").Value = 3.4591640446189009
$ws.Range("C89").Value = 3.979220649430526
$ws.Range("C90").Value = 4.4559629105553
$ws.Range("C91").Value = 4.9060512155716607
$ws.Range("C92").Value = 5.2461476225496018
$ws.Range("C93").Value = 5.2462487482324294
$ws.Range("C94").Value = 5.25635008410128
$ws.Range("C95").Value = 5.2497827595167523
$ws.Range("C96").Value = 5.0732118594682962
$ws.Range("C97").Value = 4.4966358992788003
$ws.Range("C98").Value = 3.176719879663592
$ws.Range("C99").Value = 2.0867963527678501
$ws.Range("C100").Value = 1.9401999578886997
$ws.Range("C101").Value = 0.50805703250520473
$ws.Range("C102").Value = 1.6214524346002568
$ws.Range("C103").Value = 0.21241196858492462
$ws.Range("C104").Value = -0.70387970418479551
$ws.Range("C105").Value = -0.76185294950582572
$ws.Range("C106").Value = -0.59611952736741136
$ws.Range("C107").Value = -2.0555562698919339
$ws.Range("C108").Value = -1.8806970553885805
$ws.Range("C109").Value = -2.5222233940824235
$ws.Range("C110").Value = -1.9730592759785393
$ws.Range("C111").Value = -1.5781234212484119
$ws.Range("C112").Value = -2.8532208889953004
$ws.Range("C113").Value = -2.4156176140262997
$ws.Range("C114").Value = -3.3071839649223667
$ws.Range("C115").Value = -3.0846357938296265
$ws.Range("C116").Value = -2.6387466274708449
$ws.Range("C117").Value = -3.852702029197963
$ws.Range("C118").Value = -2.3216853288092976
$ws.Range("C119").Value = -1.5982884967712607
$ws.Range("C120").Value = -1.2036652146697913
$ws.Range("C121").Value = -1.2885194897732766
$ws.Range("C122").Value = -1.3383717508449688
$ws.Range("C123").Value = -1.2696528412588282
$ws.Range("C124").Value = -0.85968761417978268
$ws.Range("C125").Value = -0.3992858350052586
$ws.Range("C126").Value = 0.1344667900119445
$ws.Range("C127").Value = 0.01834107841980348
$ws.Range("C128").Value = 0.0494541096738077
$ws.Range("C147").Value = 7.7308233574941987
$ws.Range("C148").Value = -6.9493111435087185
$ws.Range("C149").Value = -4.4095068002942162
$ws.Range("C150").Value = -3.605580720065904
$ws.Range("C151").Value = -2.7587759143008705
$ws.Range("C152").Value = -1.2762470788000191
